# Add a new translation row (TEXT ID = "SingleUseId2") to the "Translation"
# sheet, right below the existing header row.
#   B4 TEXT ID        -> SingleUseId2
#   C4 TYPOGRAPHY NAME -> Default
#   D4 ALIGNMENT       -> Right
#   E4 DIRECTION       -> LTR
#   F4 GB (text)       -> 99999

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B4").Value = "SingleUseId2"
$ws.Range("C4").Value = "Default"
$ws.Range("E4").Value = "LTR"

# Force the numeric-looking "99999" to be stored as text (matches the
# shared-string cell type in the target workbook) instead of a number.
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "99999"
$ws.Range("F4").Style = "Normal"

$ws.Range("D4").Value = "Right"
